# Hotfix: Thu Oct 17 17:53:35 RTZ 2024
#
# Adds a new "test record" row to the Python sheet, and on the Bash sheet
# rewrites the last existing row with a real bash test command plus appends
# a further "delete_bash_command" row.

$wb = $excel.ActiveWorkbook

$wsPython = $wb.Worksheets.Item("Python")
$wsBash   = $wb.Worksheets.Item("Bash")

# --- Python sheet: append row 38 (A1:D37 -> A1:D38) ---
$wsPython.Cells.Item(38, 1).Value = 2099
$wsPython.Cells.Item(38, 2).Value = "Тестовая запись"
$wsPython.Cells.Item(38, 3).Value = "Тестовая запись"
$wsPython.Cells.Item(38, 4).Value = "Тестовая запись"

# --- Bash sheet: rewrite row 70 ---
$wsBash.Cells.Item(70, 1).Value = 102
$wsBash.Cells.Item(70, 2).Value = "#!/bin/bash`nvenv/Scripts/python.exe app.py &"
$wsBash.Cells.Item(70, 3).Value = "Тестовая bash команда"

# --- Bash sheet: append row 71 (A1:C70 -> A1:C71) ---
$wsBash.Cells.Item(71, 1).Value = 105
$wsBash.Cells.Item(71, 2).Value = "delete_bash_command"
$wsBash.Cells.Item(71, 3).Value = "delete_bash_command"
